$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the cryptos list (price / 1h-volume columns, and a few reordered
# rows) to match the latest scrape. Numeric-looking price strings are
# prefixed with a leading apostrophe so Excel keeps storing them as text
# (matching the source data) instead of auto-converting them to numbers.

$ws.Range('D2').Value = '68.033.31'
$ws.Range('E2').Value = '  +1.14%  '
$ws.Range('D3').Value = '3.900.49'
$ws.Range('E3').Value = '  -1.07%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').Value = "'487.32"
$ws.Range('D6').Value = "'146.11"
$ws.Range('E6').Value = '  -0.46%  '
$ws.Range('D7').Value = "'0.619"
$ws.Range('E7').Value = '  -1.46%  '
$ws.Range('D8').Value = "'0.997"
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').Value = "'0.724"
$ws.Range('E9').Value = '  -1.51%  '
$ws.Range('D10').Value = "'0.165"
$ws.Range('E10').Value = '  -1.29%  '
$ws.Range('D11').Value = "'0.0000344"
$ws.Range('E11').Value = '  +0.58%  '
$ws.Range('D12').Value = "'42.72"
$ws.Range('E12').Value = '  -1.91%  '
$ws.Range('D13').Value = "'10.67"
$ws.Range('E13').Value = '  +2.02%  '
$ws.Range('D14').Value = '4.511.91'
$ws.Range('E14').Value = '  -0.95%  '
$ws.Range('D15').Value = '3.917.46'
$ws.Range('E15').Value = '  -1.17%  '
$ws.Range('B16').Value = 'Uniswap'
$ws.Range('C16').Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range('D16').Value = "'14.02"
$ws.Range('E16').Value = '  -7.64%  '
$ws.Range('B17').Value = 'TRON'
$ws.Range('C17').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D17').Value = "'0.136"
$ws.Range('E17').Value = '  -1.34%  '
$ws.Range('D18').Value = "'19.73"
$ws.Range('E18').Value = '  -0.86%  '
$ws.Range('E19').Value = '  -3.02%  '
$ws.Range('D20').Value = '68.159.84'
$ws.Range('E20').Value = '  +0.96%  '
$ws.Range('D21').Value = "'428.25"
$ws.Range('E21').Value = '  -2.20%  '
$ws.Range('D22').Value = "'3.54"
$ws.Range('E22').Value = '  +4.39%  '
$ws.Range('D23').Value = "'14.93"
$ws.Range('E23').Value = '  +2.85%  '
$ws.Range('D24').Value = "'87.24"
$ws.Range('E24').Value = '  -0.52%  '
$ws.Range('D25').Value = "'11.22"
$ws.Range('E25').Value = '  +14.32%  '
$ws.Range('D26').Value = "'11.25"
$ws.Range('E26').Value = '  +9.03%  '
$ws.Range('D27').Value = "'3.61"
$ws.Range('E27').Value = '  +0.41%  '
$ws.Range('D28').Value = "'38.05"
$ws.Range('E28').Value = '  -2.08%  '
$ws.Range('E29').Value = '  -0.51%  '
$ws.Range('D30').Value = "'717.43"
$ws.Range('E30').Value = '  -0.29%  '
$ws.Range('D31').Value = "'13.67"
$ws.Range('E31').Value = '  +1.01%  '
$ws.Range('D32').Value = "'0.130"
$ws.Range('E32').Value = '  -2.45%  '
$ws.Range('D33').Value = "'2.90"
$ws.Range('E33').Value = '  +2.85%  '
$ws.Range('D34').Value = "'6.30"
$ws.Range('E34').Value = '  +17.51%  '
$ws.Range('D35').Value = "'41.58"
$ws.Range('E35').Value = '  -2.97%  '
$ws.Range('D36').Value = '0.0₃0864'
$ws.Range('E36').Value = '  +7.49%  '
$ws.Range('D37').Value = "'60.28"
$ws.Range('E37').Value = '  +4.17%  '
$ws.Range('D38').Value = "'0.407"
$ws.Range('E38').Value = '  +20.49%  '
$ws.Range('D39').Value = "'0.148"
$ws.Range('E39').Value = '  -2.08%  '
$ws.Range('E41').Value = '  +15.43%  '
$ws.Range('D42').Value = "'0.0477"
$ws.Range('E42').Value = '  +0.02%  '
$ws.Range('E43').Value = '  +2.25%  '
$ws.Range('E44').Value = '  +2.66%  '
$ws.Range('B45').Value = 'Stellar'
$ws.Range('C45').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D45').Value = "'0.140"
$ws.Range('E45').Value = '  -1.35%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').Value = "'1.00"
$ws.Range('E46').Value = '  +0.16%  '
$ws.Range('D47').Value = "'3.31"
$ws.Range('E47').Value = '  +4.90%  '
$ws.Range('D48').Value = "'3.39"
$ws.Range('E48').Value = '  -5.42%  '
$ws.Range('D49').Value = "'2.12"
$ws.Range('E49').Value = '  -3.91%  '
$ws.Range('D50').Value = "'144.52"
$ws.Range('E50').Value = '  -2.32%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').Value = '0.0₆0332'
$ws.Range('E51').Value = '  +23.83%  '
